$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 0.74160000000000004
$ws.Range("H4").Value = 0.76219999999999999
$ws.Range("G5").Value = 0.69989999999999997
$ws.Range("H5").Value = 0.75829999999999997
$ws.Range("G6").Value = 0.82140000000000002
$ws.Range("H6").Value = 0.84219999999999995
$ws.Range("G7").Value = 0.82899999999999996
$ws.Range("H7").Value = 0.84909999999999997

$ws.Range("K11").Select()
